# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.134.36"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.471.14"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "583.20"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "174.45"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.89%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("E10").Value = "  +0.35%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "4.95"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("E12").Value = "  +0.69%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.925.98"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "25.42"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.72%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "67.067.01"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("E16").Value = "  -0.08%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.473.22"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("E18").Value = "  -2.18%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.47"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "349.12"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("E22").Value = "  +0.16%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "69.41"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.69%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "4.20"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.80"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.22"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.596.03"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.18%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0904"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "500.49"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -3.43%  "

$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("E34").Value = "  +0.02%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.120"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "160.80"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.14%  "

$ws.Range("E37").Value = "  +0.01%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "18.12"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.33"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("E42").Value = "  +0.01%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "4.83"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("E44").Value = "  +0.65%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "142.34"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("E48").Value = "  -0.06%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0739"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("E50").Value = "  -1.62%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.582"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
